# Marplatenses.xlsx - "Carga inicial de datos y cálculo de costo listos,
# falta capacidad de manipular bases de datos luego de cargadas"
#
# Fixes the "1/2 u" quantity typo (drop the stray space) for the two
# ingredients measured as half a unit (Limón / Naranja), leaving the
# "1 tsp" quantity (Escencia de Vainilla) untouched, and moves the
# worksheet's active selection to J6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the quantities; H2 = Limón, K2 = Naranja were both "1/2 u".
# Re-enter them without the space, matching "1/2u" everywhere else the
# fraction-style quantity is used.
$ws.Range("H2").Value = "1/2u"
$ws.Range("K2").Value = "1/2u"

# Move the active cell / selection to J6.
$ws.Range("J6").Select()
